# Updated cryptos list (price / 1h-volume-change refresh), matching the
# "Updated cryptos list ... with GitHub Actions" commit.
# For D-column values that look like plain numbers, force the cell's
# number format to Text ("@") before assigning, so Excel keeps the exact
# printed string (e.g. "561.40", "0.999") instead of silently coercing it
# to a floating-point number and losing trailing zeros / precision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.862.82'
$ws.Range("E2").Value = '  -0.74%  '

$ws.Range("D3").Value = '2.405.53'
$ws.Range("E3").Value = '  -0.78%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '561.40'
$ws.Range("E5").Value = '  +0.94%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.81'
$ws.Range("E6").Value = '  -1.28%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("E8").Value = '  -0.97%  '

$ws.Range("E9").Value = '  -0.99%  '

$ws.Range("E10").Value = '  -1.85%  '

$ws.Range("E11").Value = '  -3.00%  '

$ws.Range("E12").Value = '  -1.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '25.43'
$ws.Range("E13").Value = '  -3.20%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000171'
$ws.Range("E14").Value = '  -1.87%  '

$ws.Range("D15").Value = '2.841.11'
$ws.Range("E15").Value = '  -0.75%  '

$ws.Range("D16").Value = '61.800.52'
$ws.Range("E16").Value = '  -0.69%  '

$ws.Range("D17").Value = '2.405.70'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.19'
$ws.Range("E18").Value = '  +1.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '320.60'
$ws.Range("E19").Value = '  -1.25%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.80'
$ws.Range("E20").Value = '  +1.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.11'
$ws.Range("E21").Value = '  -2.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.61'
$ws.Range("E23").Value = '  +1.14%  '

$ws.Range("E24").Value = '  -2.45%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.65'
$ws.Range("E25").Value = '  -5.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '563.21'
$ws.Range("E26").Value = '  -1.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.07%  '

$ws.Range("D28").Value = '2.515.33'

$ws.Range("D29").Value = '0.0₃0925'
$ws.Range("E29").Value = '  -1.42%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.15'
$ws.Range("E30").Value = '  -2.89%  '

$ws.Range("E31").Value = '  -4.82%  '

$ws.Range("E32").Value = '  -0.84%  '

$ws.Range("E33").Value = '  +0.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.50'
$ws.Range("E34").Value = '  -4.08%  '

$ws.Range("E35").Value = '  +0.03%  '

$ws.Range("E36").Value = '  -1.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '152.23'
$ws.Range("E37").Value = '  +1.60%  '

$ws.Range("E38").Value = '  -5.69%  '

$ws.Range("E39").Value = '  -1.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.47'
$ws.Range("E40").Value = '  -1.57%  '

$ws.Range("E41").Value = '  -5.50%  '

$ws.Range("E42").Value = '  -0.07%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '147.72'
$ws.Range("E43").Value = '  -2.28%  '

$ws.Range("E44").Value = '  -4.13%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.58'
$ws.Range("E45").Value = '  -1.52%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0527'
$ws.Range("E46").Value = '  -2.29%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.75'
$ws.Range("E47").Value = '  -3.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.588'
$ws.Range("E48").Value = '  -0.14%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0916'
$ws.Range("E49").Value = '  +0.04%  '

$ws.Range("E50").Value = '  -1.91%  '

$ws.Range("E51").Value = '  +0.40%  '

